$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.887.91"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "3.331.15"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "556.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "3.323.89"
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.179"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.860.76"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "593.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.71%  "
$ws.Range("D17").Value = "65.916.34"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.117"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "3.319.59"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("E21").Value = "  -2.54%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.57%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("E28").Value = "  +2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "30.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("E32").Value = "  +6.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "584.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "3.708.32"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.99%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "33.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.29%  "
$ws.Range("D42").Value = "0.0₃0700"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.338"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0414"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("E51").Value = "  -0.17%  "
